# Apply cryptocurrency price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.528.42"
$ws.Range("D2").Style = $ws.Range("D4").Style
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.862.55"
$ws.Range("D3").Style = $ws.Range("D4").Style
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.93"
$ws.Range("D5").Style = $ws.Range("D4").Style
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.08"
$ws.Range("D6").Style = $ws.Range("D4").Style
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.861.40"
$ws.Range("D7").Style = $ws.Range("D4").Style
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = $ws.Range("D4").Style
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  +3.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000291"
$ws.Range("D12").Style = $ws.Range("D4").Style
$ws.Range("E12").Value = "  +16.97%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.31"
$ws.Range("D14").Style = $ws.Range("D4").Style
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.512.65"
$ws.Range("D15").Style = $ws.Range("D4").Style
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.869.50"
$ws.Range("D16").Style = $ws.Range("D4").Style
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.604.23"
$ws.Range("D17").Style = $ws.Range("D4").Style
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.46"
$ws.Range("D18").Style = $ws.Range("D4").Style
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("D19").Style = $ws.Range("D4").Style
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.92"
$ws.Range("D21").Style = $ws.Range("D4").Style
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.44"
$ws.Range("D22").Style = $ws.Range("D4").Style
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.736"
$ws.Range("D23").Style = $ws.Range("D4").Style
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000162"
$ws.Range("D24").Style = $ws.Range("D4").Style
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.98"
$ws.Range("D25").Style = $ws.Range("D4").Style
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.27"
$ws.Range("D27").Style = $ws.Range("D4").Style
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.59"
$ws.Range("D28").Style = $ws.Range("D4").Style
$ws.Range("E28").Value = "  +5.85%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.015.62"
$ws.Range("D31").Style = $ws.Range("D4").Style
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.28"
$ws.Range("D34").Style = $ws.Range("D4").Style
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.48"
$ws.Range("D35").Style = $ws.Range("D4").Style
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.830.48"
$ws.Range("D36").Style = $ws.Range("D4").Style
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.97"
$ws.Range("D37").Style = $ws.Range("D4").Style
$ws.Range("E37").Value = "  +21.33%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.02"
$ws.Range("D39").Style = $ws.Range("D4").Style
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.00"
$ws.Range("D40").Style = $ws.Range("D4").Style
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.322"
$ws.Range("D43").Style = $ws.Range("D4").Style
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("D44").Style = $ws.Range("D4").Style
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("B45").Value = "FLOKI"
$ws.Range("C45").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000299"
$ws.Range("D45").Style = $ws.Range("D4").Style
$ws.Range("E45").Value = "  +10.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.82"
$ws.Range("D46").Style = $ws.Range("D4").Style
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = $ws.Range("D4").Style
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "422.89"
$ws.Range("D48").Style = $ws.Range("D4").Style
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.74"
$ws.Range("D49").Style = $ws.Range("D4").Style
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0361"
$ws.Range("D50").Style = $ws.Range("D4").Style
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.43"
$ws.Range("D51").Style = $ws.Range("D4").Style
$ws.Range("E51").Value = "  -0.53%  "
